$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range("A100").NumberFormat = "@"
    $ws.Range("A100").Value = $val
    $ws.Range("A100").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("A100").Delete()
}

$ws.Range('E2').Value = '2026-02-09 07:18:43'
$ws.Range('E3').Value = '2026-02-09 07:18:46'
$ws.Range('M3').Value = '-5.2 °C 6:58 TU'
$ws.Range('E4').Value = '2026-02-09 07:18:48'
$ws.Range('E5').Value = '2026-02-09 07:18:51'
$ws.Range('O5').Value = '-5.2 °C'
$ws.Range('E6').Value = '2026-02-09 07:18:54'
$ws.Range('J6').Value = '1007.7 hPa'
$ws.Range('N6').Value = '4.4 °C 6:59 TU'
$ws.Range('O6').Value = '6.1 °C'
$ws.Range('E7').Value = '2026-02-09 07:18:56'
Set-TextValue 'H7' '68%'
$ws.Range('L7').Value = '14.8 km/h - 302º 6:47 TU'
$ws.Range('N7').Value = '10.7 °C 6:40 TU'
$ws.Range('E8').Value = '2026-02-09 07:18:59'
Set-TextValue 'H8' '77%'
$ws.Range('J8').Value = '1008.4 hPa'
$ws.Range('N8').Value = '6.2 °C 6:57 TU'
$ws.Range('E9').Value = '2026-02-09 07:19:01'
$ws.Range('N9').Value = '2.3 °C 6:59 TU'
$ws.Range('O9').Value = '5.8 °C'
$ws.Range('E10').Value = '2026-02-09 07:19:04'
Set-TextValue 'H10' '96%'
$ws.Range('O10').Value = '4.9 °C'
$ws.Range('E11').Value = '2026-02-09 07:19:07'
$ws.Range('E12').Value = '2026-02-09 07:19:09'
$ws.Range('N12').Value = '2.7 °C 6:39 TU'
$ws.Range('O12').Value = '6.5 °C'
$ws.Range('E13').Value = '2026-02-09 07:19:12'
$ws.Range('J13').Value = '1011.1 hPa'
$ws.Range('E14').Value = '2026-02-09 07:19:14'
$ws.Range('E15').Value = '2026-02-09 07:19:17'
$ws.Range('N15').Value = '2.3 °C 6:39 TU'
$ws.Range('O15').Value = '4.8 °C'
$ws.Range('E16').Value = '2026-02-09 07:19:19'
Set-TextValue 'H16' '66%'
$ws.Range('N16').Value = '-6.7 °C 6:54 TU'
$ws.Range('E17').Value = '2026-02-09 07:19:22'
$ws.Range('E18').Value = '2026-02-09 07:19:25'
Set-TextValue 'H18' '99%'
$ws.Range('O18').Value = '5.6 °C'
$ws.Range('E19').Value = '2026-02-09 07:19:27'
$ws.Range('N19').Value = '2.1 °C 6:49 TU'
$ws.Range('E20').Value = '2026-02-09 07:19:30'
$ws.Range('O20').Value = '-6.1 °C'
$ws.Range('E21').Value = '2026-02-09 07:19:32'
Set-TextValue 'H21' '95%'
$ws.Range('J21').Value = '1009.9 hPa'
$ws.Range('E22').Value = '2026-02-09 07:19:34'
$ws.Range('E23').Value = '2026-02-09 07:19:36'
Set-TextValue 'H23' '85%'
$ws.Range('E24').Value = '2026-02-09 07:19:39'
$ws.Range('E25').Value = '2026-02-09 07:19:41'
Set-TextValue 'H25' '74%'
$ws.Range('M25').Value = '-2.9 °C 6:55 TU'
$ws.Range('O25').Value = '-4.1 °C'
$ws.Range('E26').Value = '2026-02-09 07:19:44'
Set-TextValue 'H26' '92%'
$ws.Range('E27').Value = '2026-02-09 07:19:47'
Set-TextValue 'H27' '84%'
$ws.Range('O27').Value = '-4.1 °C'
$ws.Range('E28').Value = '2026-02-09 07:19:50'
$ws.Range('O28').Value = '3.3 °C'
$ws.Range('E29').Value = '2026-02-09 07:19:53'
$ws.Range('E30').Value = '2026-02-09 07:19:55'
$ws.Range('N30').Value = '4.0 °C 6:46 TU'
$ws.Range('O30').Value = '6.1 °C'
$ws.Range('E31').Value = '2026-02-09 07:19:58'
Set-TextValue 'H31' '73%'
$ws.Range('E32').Value = '2026-02-09 07:20:01'
$ws.Range('E33').Value = '2026-02-09 07:20:04'
$ws.Range('E34').Value = '2026-02-09 07:20:06'
Set-TextValue 'H34' '74%'
$ws.Range('O34').Value = '-3.5 °C'
$ws.Range('E35').Value = '2026-02-09 07:20:09'
Set-TextValue 'H35' '68%'
$ws.Range('I35').Value = '0.1 mm'
$ws.Range('N35').Value = '2.3 °C 6:59 TU'
$ws.Range('O35').Value = '3.5 °C'
$ws.Range('E36').Value = '2026-02-09 07:20:12'
Set-TextValue 'H36' '87%'
$ws.Range('N36').Value = '4.1 °C 6:40 TU'
$ws.Range('O36').Value = '7.6 °C'
$ws.Range('E37').Value = '2026-02-09 07:20:15'
$ws.Range('O37').Value = '2.8 °C'
$ws.Range('E38').Value = '2026-02-09 07:20:18'
$ws.Range('N38').Value = '3.1 °C 6:40 TU'
$ws.Range('O38').Value = '5.4 °C'
$ws.Range('E39').Value = '2026-02-09 07:20:20'
$ws.Range('E40').Value = '2026-02-09 07:20:23'
$ws.Range('J40').Value = '1010.6 hPa'
$ws.Range('E41').Value = '2026-02-09 07:20:26'
$ws.Range('E42').Value = '2026-02-09 07:20:29'
$ws.Range('N42').Value = '3.1 °C 6:59 TU'
$ws.Range('O42').Value = '5.8 °C'
$ws.Range('E43').Value = '2026-02-09 07:20:31'
$ws.Range('N43').Value = '5.5 °C 6:30 TU'
$ws.Range('O43').Value = '6.2 °C'
$ws.Range('E44').Value = '2026-02-09 07:20:34'
$ws.Range('O44').Value = '-6.9 °C'
$ws.Range('E45').Value = '2026-02-09 07:20:37'
$ws.Range('J45').Value = '1009.3 hPa'
$ws.Range('O45').Value = '0.3 °C'
$ws.Range('E46').Value = '2026-02-09 07:20:40'
